# Add a new "TextBox 1" shape to slide 1:
#   "К 29.04 будет готова"
# at off(x=4563611, y=2491530) ext(cx=2169697, cy=369332) EMU.
#
# PowerPoint COM positions/sizes are expressed in points (1 pt = 12700 EMU).
# The target EMU values below are reproduced from the full-precision
# points equivalents (target_emu / 12700) so that Left/Top/Width/Height
# round-trip back to exactly those EMU numbers.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)   # 1 = msoTextOrientationHorizontal
$shp.Name = "TextBox 1"

$shp.Left   = 359.3394623188977     # -> 4563611 EMU
$shp.Top    = 196.18346456692913    # -> 2491530 EMU
$shp.Width  = 170.84228346456692    # -> 2169697 EMU
$shp.Height = 29.081259842519685    # -> 369332 EMU

# bodyPr: wrap="none" + <a:spAutoFit/>, spPr: <a:noFill/>
$shp.TextFrame.WordWrap = 0         # msoFalse -> wrap="none"
$shp.TextFrame.AutoSize = 1         # ppAutoSizeShapeToFitText -> <a:spAutoFit/>
$shp.Fill.Visible = 0               # msoFalse -> <a:noFill/>

$shp.TextFrame.TextRange.Text = "К 29.04 будет готова"
$shp.TextFrame.TextRange.LanguageID = "ru-RU"   # rPr lang="ru-RU"
